# Updates cryptos list: Price (D) and Volume(1h) (E) columns, plus row 12/13 coin swap (Chainlink <-> TRON)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.036.85'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '2.839.31'
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '361.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.571'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.36%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.602'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.61'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0861'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.131'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").Value = '3.287.25'
$ws.Range("E15").Value = '  +1.79%  '
$ws.Range("D16").Value = '2.825.49'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").Value = '51.856.91'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.63%  '
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.00%  '
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.91%  '
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.62%  '
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '53.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.45%  '
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.02'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("E33").Value = '  +3.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0449'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +20.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.32'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0841'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.24%  '
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("E39").Value = '  -2.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("E42").Value = '  +1.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '127.78'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.55'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.96%  '
$ws.Range("E45").Value = '  -3.47%  '
$ws.Range("D46").Value = '2.122.35'
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("E48").Value = '  +1.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.84'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.01%  '
$ws.Range("E51").Value = '  +0.95%  '
